# Apply the "output generated at 456a3b4" data refresh to all sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value  = 854
$ws.Range("F3").Value  = 13742
$ws.Range("F4").Value  = 13533
$ws.Range("F7").Value  = 39
$ws.Range("F8").Value  = 592
$ws.Range("F9").Value  = 80
$ws.Range("F11").Value = 50
$ws.Range("F12").Value = 753
$ws.Range("F13").Value = 2141
$ws.Range("F14").Value = 88
$ws.Range("G16").Value = 29.9
$ws.Range("I17").Value = "//i1.hdslb.com/bfs/openplatform/202406/sk6wpxN91717486689960.jpeg"
$ws.Range("F19").Value = 516
$ws.Range("F21").Value = 386
$ws.Range("F22").Value = 319
$ws.Range("F24").Value = 826
$ws.Range("F25").Value = 76

# --- Sheet "演出" (Performances) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value  = 162
$ws.Range("F7").Value  = 1472
$ws.Range("F11").Value = 63

# --- Sheet "本地生活" (Local life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 219
$ws.Range("F3").Value = 103

# --- Sheet "全部类型" (All types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value  = 219
$ws.Range("F3").Value  = 854
$ws.Range("F4").Value  = 13742
$ws.Range("F5").Value  = 13533
$ws.Range("F8").Value  = 39
$ws.Range("F9").Value  = 592
$ws.Range("F10").Value = 80
$ws.Range("F12").Value = 50
$ws.Range("F13").Value = 753
$ws.Range("F16").Value = 2141
$ws.Range("F17").Value = 88
$ws.Range("G19").Value = 29.9
$ws.Range("I20").Value = "//i1.hdslb.com/bfs/openplatform/202406/sk6wpxN91717486689960.jpeg"
$ws.Range("F24").Value = 103
$ws.Range("F25").Value = 103
$ws.Range("F26").Value = 516
$ws.Range("F28").Value = 386
$ws.Range("F29").Value = 319
$ws.Range("F31").Value = 826
$ws.Range("F32").Value = 162
$ws.Range("F33").Value = 1472
$ws.Range("F37").Value = 76
$ws.Range("F38").Value = 63
